# ============================================================================
# Hormoz Steel quarterly income statement (rial, cumulative) - refresh update
# Two additional historical report columns are prepended (periods ending
# 1399/06 through 1400/06), the previously-latest five periods shift right,
# and the "EPS on latest capital" row is recomputed with the new algorithm.
# ============================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1. Row 8 (period label) / Row 9 (publish date) headers for the I:M block
#    (the figures that used to sit in D:H). These are written FIRST, as
#    plain text, before the format-copy below - a couple of the publish-
#    date strings are bare "yyyy-mm-dd" text that Excel would otherwise
#    auto-coerce to a date serial, so they're entered with a leading quote
#    to force text, same as typing them in by hand would require.
# ----------------------------------------------------------------------
$row8_IM = @{
    "I" = "9 ماهه منتهی به 1400/09"
    "J" = "12 ماهه منتهی به 1400/12"
    "K" = "3 ماهه منتهی به 1401/03"
    "L" = "6 ماهه منتهی به 1401/06"
    "M" = "9 ماهه منتهی به 1401/09"
}
$row9_IM = @{
    "I" = "1401-10-28 (2)"
    "J" = "1401-10-28 (6)"
    "K" = "'1401-04-29"
    "L" = "1401-09-13 (2)"
    "M" = "'1401-10-28"
}
foreach ($c in $row8_IM.Keys) { $ws.Range($c + "8").Value = $row8_IM[$c] }
foreach ($c in $row9_IM.Keys) { $ws.Range($c + "9").Value = $row9_IM[$c] }

# ----------------------------------------------------------------------
# 2. Copy the existing D:H formatting onto the new I:M block so the
#    columns that just received the previous figures (shifted right by
#    five columns) keep identical styling (fills/borders/alignment/etc),
#    and so the quote-prefixed cells above end up on the normal style
#    instead of a "quote prefix" variant.
# ----------------------------------------------------------------------
$ws.Range("D1:H28").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# 3. Column widths: the 28/28/29/28/28-char width pattern now repeats
#    twice (new D:H block, then the shifted I:M block). Range.ColumnWidth
#    adds an implicit 0.8333333... pad versus the raw XML "width", so we
#    subtract it up front to land on the exact target width.
# ----------------------------------------------------------------------
$pad = 0.8333333333333334
$colWidths = @{ "D" = 28; "E" = 28; "F" = 29; "G" = 28; "H" = 28; "I" = 28; "J" = 29; "K" = 28; "L" = 28; "M" = 28 }
foreach ($c in $colWidths.Keys) {
    $ws.Columns($c).ColumnWidth = $colWidths[$c] - $pad
}

# ----------------------------------------------------------------------
# 4. Row 8 / Row 9 headers for the new D:H block (newly-disclosed older
#    periods). None of these strings are bare "yyyy-mm-dd" dates, so no
#    quote-prefix is required.
# ----------------------------------------------------------------------
$row8_DH = @{
    "D" = "6 ماهه منتهی به 1399/06"
    "E" = "9 ماهه منتهی به 1399/09"
    "F" = "12 ماهه منتهی به 1399/12"
    "G" = "3 ماهه منتهی به 1400/03"
    "H" = "6 ماهه منتهی به 1400/06"
}
$row9_DH = @{
    "D" = "1400-09-14 (5)"
    "E" = "1400-10-29 (2)"
    "F" = "1401-03-15 (9)"
    "G" = "1401-04-29 (2)"
    "H" = "1401-09-13 (5)"
}
foreach ($c in $row8_DH.Keys) { $ws.Range($c + "8").Value = $row8_DH[$c] }
foreach ($c in $row9_DH.Keys) { $ws.Range($c + "9").Value = $row9_DH[$c] }

# ----------------------------------------------------------------------
# 5. Data rows 11-27.
#    D:H  -> newly-disclosed older periods (1399/06 .. 1400/06)
#    I:M  -> figures that used to live in D:H, shifted right - EXCEPT row
#            27 (EPS on latest capital), which was recomputed under the
#            new read_price algorithm and is NOT the old D:H data.
# ----------------------------------------------------------------------
$newPeriodData = @{
    11 = @(42850994, 76745696, 134522191, 52463805, 102598264)
    12 = @(-24923113, -43530550, -76944388, -32047835, -51859003)
    13 = @(17927881, 33215146, 57577803, 20415970, 50739261)
    14 = @(-873287, -1508367, -2341745, -954441, -2018289)
    15 = @(0, 0, 0, 0, 0)
    16 = @(6905309, 7018729, 4641506, 135721, -1721391)
    17 = @(23959903, 38725508, 59877564, 19597250, 46999581)
    18 = @(-819643, -1517191, -2384228, -771051, -906433)
    19 = @(353824, 811046, -315283, -54433, -43198)
    20 = @(23494084, 38019363, 57178053, 18771766, 46049950)
    21 = @(-2307568, -4268196, -5153126, -2162963, -4436910)
    22 = @(21186516, 33751167, 52024927, 16608803, 41613040)
    23 = @(0, 0, 0, 0, 0)
    24 = @(21186516, 33751167, 52024927, 16608803, 41613040)
    25 = @(847, 1350, 2081, 664, 1665)
    26 = @(25000000, 25000000, 25000000, 25000000, 25000000)
    27 = @(432, 689, 1062, 339, 849)
}
$shiftedData = @{
    11 = @(163865794, 221404720, 48797480, 89964578, 141519631)
    12 = @(-87980472, -127265274, -32448257, -62418157, -107515985)
    13 = @(75885322, 94139446, 16349223, 27546421, 34003646)
    14 = @(-3045695, -4748509, -1298586, -3035681, -4843415)
    15 = @(0, 0, 0, 0, 0)
    16 = @(-1184597, 407127, -388640, 1707857, 7298428)
    17 = @(71655030, 89798064, 14661997, 26218597, 36458659)
    18 = @(-1051817, -1753460, -828780, -1757092, -3121481)
    19 = @(150461, -661091, -130968, 3226834, 1773359)
    20 = @(70753674, 87383513, 13702249, 27688339, 35110537)
    21 = @(-8071151, -7146516, -1620328, -2623456, -4004399)
    22 = @(62682523, 80236997, 12081921, 25064883, 31106138)
    23 = @(0, 0, 0, 0, 0)
    24 = @(62682523, 80236997, 12081921, 25064883, 31106138)
    25 = @(2507, 3209, 483, 1003, 1244)
    26 = @(25000000, 25000000, 25000000, 25000000, 25000000)
    27 = @(1279, 1637, 247, 512, 635)
}
$cols5 = @("D","E","F","G","H")
$cols5b = @("I","J","K","L","M")
foreach ($r in $newPeriodData.Keys) {
    $vals = $newPeriodData[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($cols5[$i] + $r).Value = $vals[$i]
    }
}
foreach ($r in $shiftedData.Keys) {
    $vals = $shiftedData[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($cols5b[$i] + $r).Value = $vals[$i]
    }
}

